$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, shifting existing rows 100-124 down to 101-125.
$ws.Rows("100:100").Insert()

# Populate the newly inserted row 100 with the new weekly record.
$ws.Range("A100").Value = 4
$ws.Range("B100").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C100").Value = "Los Lagos"
$ws.Range("D100").Value = 44489
$ws.Range("E100").Value = 10
$ws.Range("F100").Value = 100112028
$ws.Range("G100").Value = "Sandia"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 150
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = 1000
$ws.Range("N100").Value = "$/kilo (volumen en unidades)"
$ws.Range("O100").Value = "Perú"
$ws.Range("P100").Value = 1000
$ws.Range("Q100").Value = 1
$ws.Range("R100").Value = "Hortaliza"
